# Edit script: "trouble shooting migration.py with polars"
#
# Summary of the change being applied:
#   1. The first slide (the "script flow" diagram) is duplicated. One copy keeps
#      the first half of the pipeline (Group 5, Group 6, Connector 10, Group 11,
#      Group 14) and becomes the new first slide. The other copy keeps the
#      second half (Group 6, Group 20, Connector 23, Group 27) and is moved to
#      become the new last slide.
#   2. The middle slide gets a small text edit: "by_age_group_" is removed from
#      one of the rectangle labels.

$p = $ppt.ActivePresentation

# --- Step 1: duplicate the first slide -------------------------------------
$original = $p.Slides.Item(1)
$dupRange = $original.Duplicate()
$duplicate = $dupRange.Item(1)

# After Duplicate(), slide order is: [original(256), duplicate(258), slide2(257)]

# --- Step 2: trim the duplicate (this becomes the new FIRST slide) ---------
# Keep: Group 5, Group 6, Straight Arrow Connector 10, Group 11, Group 14
# Remove: Group 17, Group 20, Straight Arrow Connector 23, Group 27
$removeFromDuplicate = @(
    "Group 17",
    "Group 20",
    "Straight Arrow Connector 23",
    "Group 27"
)
for ($i = $duplicate.Shapes.Count; $i -ge 1; $i--) {
    $shp = $duplicate.Shapes.Item($i)
    if ($removeFromDuplicate -contains $shp.Name) {
        $shp.Delete()
    }
}

# --- Step 3: trim the original (this becomes the new LAST slide) -----------
# Keep: Group 6, Group 20, Straight Arrow Connector 23, Group 27
# Remove: Group 5, Straight Arrow Connector 10, Group 11, Group 14, Group 17
$removeFromOriginal = @(
    "Group 5",
    "Straight Arrow Connector 10",
    "Group 11",
    "Group 14",
    "Group 17"
)
for ($i = $original.Shapes.Count; $i -ge 1; $i--) {
    $shp = $original.Shapes.Item($i)
    if ($removeFromOriginal -contains $shp.Name) {
        $shp.Delete()
    }
}

# --- Step 4: reorder so the trimmed original slide goes to the end ---------
$original.MoveTo($p.Slides.Count)

# Final order is now: [duplicate(258), slide2(257), original(256)]

# --- Step 5: text edit on the middle slide ----------------------------------
# Only the "Rectangle 5" shape inside "Group 3" changes -- the very similar
# looking "Rectangle 40" inside "Group 38" keeps its original text.
$midSlide = $p.Slides.Item(2)
$targetGroup = $null
foreach ($shp in $midSlide.Shapes) {
    if ($shp.Name -eq "Group 3") {
        $targetGroup = $shp
    }
}
if ($targetGroup -ne $null) {
    foreach ($item in $targetGroup.GroupItems) {
        if ($item.Name -eq "Rectangle 5") {
            $item.TextFrame.TextRange.Text = "acs_immigration_cohort_fractions_2006-2015"
        }
    }
}
